$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 4 new columns (K:N) for Northern Lat, Southern Lat, Western Lon, Eastern Lon
#    This shifts old K..R (Site ID..Contributor) right to O..V
$ws.Range("K1:N1").EntireColumn.Insert()

# 2. New header row values for inserted columns
$ws.Range("K1").Value = "Northern Lat"
$ws.Range("L1").Value = "Southern Lat"
$ws.Range("M1").Value = "Western Lon"
$ws.Range("N1").Value = "Eastern Lon"

# 3. New data row values for inserted columns
$ws.Range("K2").Value = 42.513333
$ws.Range("L2").Value = 42.513333
$ws.Range("M2").Value = -72.2183333
$ws.Range("N2").Value = -72.2183333

# 4. Update Scan resolution value (G2) 3600 -> 3200
$ws.Range("G2").Value = 3200

# 5. Replace formula-driven booleans with literal boolean values
$ws.Range("H2").Value = $false
$ws.Range("I2").Value = $true

# 6. Fix date number format to lowercase pattern
$ws.Range("D2").NumberFormat = "yyyy\-mm\-dd"

# 7. Row 2 height back to default 12.8 (remove custom height flag)
$ws.Rows.Item(2).AutoFit()

# 8. Update sheet view selection to match target (I2)
$ws.Range("I2").Select()
